# Auto-generated edit script: update computed market-price / profit columns
# across multiple sheets, matching the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1614.7333
$ws.Range("I100").Value = 1573.0714
$ws.Range("K100").Value = 1573.0714
$ws.Range("M100").Value = -1032.0714
$ws.Range("H106").Value = 3442.6365
$ws.Range("J106").Value = 2995
$ws.Range("L106").Value = 2995
$ws.Range("N106").Value = -4257
$ws.Range("H121").Value = 1770.9333
$ws.Range("I121").Value = 1148.5
$ws.Range("J121").Value = 1866.6923
$ws.Range("K121").Value = 3445.5
$ws.Range("L121").Value = 5600.0769
$ws.Range("M121").Value = -1698.5
$ws.Range("N121").Value = -9094.0769
$ws.Range("H138").Value = 5287.615
$ws.Range("I138").Value = 11254.4
$ws.Range("K138").Value = 33763.2
$ws.Range("M138").Value = -28623.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 6465.9653
$ws.Range("I61").Value = 7396.375
$ws.Range("K61").Value = 7396.375
$ws.Range("M61").Value = -7184.375
$ws.Range("H97").Value = 563.36365
$ws.Range("I97").Value = 599.7
$ws.Range("K97").Value = 599.7
$ws.Range("M97").Value = -103.7
$ws.Range("H102").Value = 3463.6667
$ws.Range("I102").Value = 2396.625
$ws.Range("K102").Value = 2396.625
$ws.Range("M102").Value = -774.625
$ws.Range("H132").Value = 1148.6
$ws.Range("I132").Value = 1052.9429
$ws.Range("K132").Value = 3158.8287
$ws.Range("M132").Value = -628.8287
$ws.Range("H136").Value = 6465.9653
$ws.Range("I136").Value = 7396.375
$ws.Range("K136").Value = 22189.125
$ws.Range("M136").Value = -19639.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 439.4
$ws.Range("J64").Value = 600
$ws.Range("L64").Value = 600
$ws.Range("N64").Value = -1050
$ws.Range("H67").Value = 439.4
$ws.Range("J67").Value = 600
$ws.Range("L67").Value = 600
$ws.Range("N67").Value = -2160
$ws.Range("H94").Value = 784.2143
$ws.Range("I94").Value = 798.3333
$ws.Range("K94").Value = 798.3333
$ws.Range("M94").Value = -347.3333
$ws.Range("H99").Value = 5217.294
$ws.Range("I99").Value = 4742.857
$ws.Range("K99").Value = 4742.857
$ws.Range("M99").Value = -3244.857
$ws.Range("H100").Value = 17992
$ws.Range("J100").Value = 17992
$ws.Range("L100").Value = 17992
$ws.Range("N100").Value = -20156
$ws.Range("H134").Value = 2100.0476
$ws.Range("I134").Value = 2096.6667
$ws.Range("K134").Value = 6290.000100000001
$ws.Range("M134").Value = -3755.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38466840
$ws.Range("I31").Value = 111115800
$ws.Range("K31").Value = 111115800
$ws.Range("M31").Value = -111115505
$ws.Range("H34").Value = 38466840
$ws.Range("I34").Value = 111115800
$ws.Range("K34").Value = 111115800
$ws.Range("M34").Value = -111115598
$ws.Range("H141").Value = 190416.84
$ws.Range("J141").Value = 209672.53
$ws.Range("L141").Value = 209672.53
$ws.Range("N141").Value = -220032.53

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 12600951
$ws.Range("I11").Value = 14001023
$ws.Range("K11").Value = 42003069
$ws.Range("M11").Value = -42002929
$ws.Range("H98").Value = 561
$ws.Range("I98").Value = 708.6
$ws.Range("J98").Value = 508.2857
$ws.Range("K98").Value = 2125.8
$ws.Range("L98").Value = 1524.8571
$ws.Range("M98").Value = -627.8000000000002
$ws.Range("N98").Value = -4520.8571
$ws.Range("H107").Value = 764.7368
$ws.Range("I107").Value = 432.7
$ws.Range("J107").Value = 1133.6666
$ws.Range("K107").Value = 1298.1
$ws.Range("L107").Value = 3400.9998
$ws.Range("M107").Value = 621.9000000000001
$ws.Range("N107").Value = -7240.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3916.4443
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 3916.4443
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H97").Value = 900.4167
$ws.Range("J97").Value = 959.25
$ws.Range("L97").Value = 959.25
$ws.Range("N97").Value = -1951.25
$ws.Range("H102").Value = 3824.7585
$ws.Range("I102").Value = 2976.76
$ws.Range("J102").Value = 9124.75
$ws.Range("K102").Value = 2976.76
$ws.Range("L102").Value = 9124.75
$ws.Range("M102").Value = -1354.76
$ws.Range("N102").Value = -12368.75
$ws.Range("H132").Value = 3885
$ws.Range("I132").Value = 2341.8572
$ws.Range("K132").Value = 7025.571599999999
$ws.Range("M132").Value = -4495.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 361.2
$ws.Range("I16").Value = 373.55554
$ws.Range("K16").Value = 373.55554
$ws.Range("M16").Value = -203.55554
$ws.Range("H41").Value = 666666
$ws.Range("I41").Value = 666666
$ws.Range("K41").Value = 666666
$ws.Range("M41").Value = -666228
$ws.Range("H46").Value = 5359.8535
$ws.Range("I46").Value = 6467.032
$ws.Range("K46").Value = 6467.032
$ws.Range("M46").Value = -6279.032
$ws.Range("H68").Value = 4961.154
$ws.Range("J68").Value = 6227.7144
$ws.Range("L68").Value = 6227.7144
$ws.Range("N68").Value = -7725.7144
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H71").Value = 4961.154
$ws.Range("J71").Value = 6227.7144
$ws.Range("L71").Value = 31138.572
$ws.Range("N71").Value = -38626.572
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H100").Value = 5772.7393
$ws.Range("I100").Value = 5084.154
$ws.Range("K100").Value = 5084.154
$ws.Range("M100").Value = -4543.154
$ws.Range("H122").Value = 1950.2858
$ws.Range("I122").Value = 1950.2858
$ws.Range("K122").Value = 5850.857400000001
$ws.Range("M122").Value = -3400.857400000001
$ws.Range("H132").Value = 26875.5
$ws.Range("I132").Value = 27275.25
$ws.Range("K132").Value = 81825.75
$ws.Range("M132").Value = -79295.75
$ws.Range("H136").Value = 58828370
$ws.Range("I136").Value = 50005644
$ws.Range("K136").Value = 150016932
$ws.Range("M136").Value = -150014382

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3065.5
$ws.Range("I96").Value = 3491
$ws.Range("J96").Value = 2923.6667
$ws.Range("K96").Value = 3491
$ws.Range("L96").Value = 2923.6667
$ws.Range("M96").Value = -2118
$ws.Range("N96").Value = -5669.6667
$ws.Range("H100").Value = 1079.4445
$ws.Range("I100").Value = 966.05884
$ws.Range("K100").Value = 1932.11768
$ws.Range("M100").Value = -1391.11768
$ws.Range("H116").Value = 340000
$ws.Range("J116").Value = 340000
$ws.Range("L116").Value = 340000
$ws.Range("N116").Value = -349178
$ws.Range("H132").Value = 3002.6064
$ws.Range("I132").Value = 2680.0728
$ws.Range("K132").Value = 8040.2184
$ws.Range("M132").Value = -5510.2184
$ws.Range("H136").Value = 5817.909
$ws.Range("I136").Value = 4597.189
$ws.Range("J136").Value = 8327.166999999999
$ws.Range("K136").Value = 13791.567
$ws.Range("L136").Value = 24981.501
$ws.Range("M136").Value = -11241.567
$ws.Range("N136").Value = -30081.501
